# Update the "想去人数" (column F) values on the 展览 / 演出 / 全部类型 sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 13047
$wsExhibit.Range("F5").Value = 86
$wsExhibit.Range("F8").Value = 26
$wsExhibit.Range("F10").Value = 13013
$wsExhibit.Range("F11").Value = 295
$wsExhibit.Range("F12").Value = 548
$wsExhibit.Range("F13").Value = 8732
$wsExhibit.Range("F14").Value = 7751
$wsExhibit.Range("F16").Value = 131
$wsExhibit.Range("F23").Value = 187

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 5

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 13047
$wsAll.Range("F6").Value = 86
$wsAll.Range("F9").Value = 26
$wsAll.Range("F11").Value = 13013
$wsAll.Range("F12").Value = 295
$wsAll.Range("F13").Value = 548
$wsAll.Range("F14").Value = 8732
$wsAll.Range("F15").Value = 7751
$wsAll.Range("F17").Value = 131
$wsAll.Range("F23").Value = 5
$wsAll.Range("F26").Value = 187
